$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "94.824.38"
Set-TextValue $ws.Range("E2") "  -2.02%  "
Set-TextValue $ws.Range("D3") "3.555.49"
Set-TextValue $ws.Range("E3") "  -1.70%  "
Set-TextValue $ws.Range("E4") "  -0.01%  "
Set-TextValue $ws.Range("D5") "235.51"
Set-TextValue $ws.Range("E5") "  -2.63%  "
Set-TextValue $ws.Range("D6") "652.69"
Set-TextValue $ws.Range("E6") "  +1.64%  "
Set-TextValue $ws.Range("E7") "  -1.89%  "
Set-TextValue $ws.Range("E8") "  -2.09%  "
Set-TextValue $ws.Range("E9") "  +0.13%  "
Set-TextValue $ws.Range("D10") "0.999"
Set-TextValue $ws.Range("E10") "  -2.57%  "
Set-TextValue $ws.Range("D11") "3.551.74"
Set-TextValue $ws.Range("E11") "  -1.62%  "
Set-TextValue $ws.Range("E12") "  +0.37%  "
Set-TextValue $ws.Range("D13") "42.13"
Set-TextValue $ws.Range("E13") "  -3.57%  "
Set-TextValue $ws.Range("D14") "6.43"
Set-TextValue $ws.Range("E14") "  +0.38%  "
Set-TextValue $ws.Range("D15") "4.220.39"
Set-TextValue $ws.Range("E15") "  -1.94%  "
Set-TextValue $ws.Range("D16") "94.807.15"
Set-TextValue $ws.Range("E16") "  -1.99%  "
Set-TextValue $ws.Range("E17") "  -1.30%  "
Set-TextValue $ws.Range("D18") "3.558.67"
Set-TextValue $ws.Range("E18") "  -1.55%  "
Set-TextValue $ws.Range("D19") "7.90"
Set-TextValue $ws.Range("E19") "  -0.81%  "
Set-TextValue $ws.Range("E20") "  -3.46%  "
Set-TextValue $ws.Range("E21") "  -3.32%  "
Set-TextValue $ws.Range("E22") "  -1.48%  "
Set-TextValue $ws.Range("D23") "506.10"
Set-TextValue $ws.Range("E23") "  -2.16%  "
Set-TextValue $ws.Range("D24") "0.478"
Set-TextValue $ws.Range("E24") "  -4.74%  "
Set-TextValue $ws.Range("D25") "6.78"
Set-TextValue $ws.Range("E25") "  -0.26%  "
Set-TextValue $ws.Range("E26") "  -3.19%  "
Set-TextValue $ws.Range("D27") "94.72"
Set-TextValue $ws.Range("E27") "  -3.79%  "
Set-TextValue $ws.Range("D28") "12.52"
Set-TextValue $ws.Range("E28") "  -0.30%  "
Set-TextValue $ws.Range("D29") "3.748.36"
Set-TextValue $ws.Range("E29") "  -1.53%  "
Set-TextValue $ws.Range("E30") "  -5.34%  "
Set-TextValue $ws.Range("D31") "0.143"
Set-TextValue $ws.Range("E31") "  -1.74%  "
Set-TextValue $ws.Range("D32") "11.44"
Set-TextValue $ws.Range("E32") "  -2.11%  "
Set-TextValue $ws.Range("E33") "  +0.10%  "
Set-TextValue $ws.Range("E34") "  -0.40%  "
Set-TextValue $ws.Range("D35") "0.175"
Set-TextValue $ws.Range("E35") "  -4.85%  "
Set-TextValue $ws.Range("D36") "31.64"
Set-TextValue $ws.Range("E36") "  +3.55%  "
Set-TextValue $ws.Range("D37") "1.68"
Set-TextValue $ws.Range("E37") "  +13.19%  "
Set-TextValue $ws.Range("D38") "0.552"
Set-TextValue $ws.Range("E38") "  -3.53%  "
Set-TextValue $ws.Range("D39") "8.43"
Set-TextValue $ws.Range("E39") "  +6.30%  "
Set-TextValue $ws.Range("D40") "579.43"
Set-TextValue $ws.Range("E40") "  -0.61%  "
Set-TextValue $ws.Range("E41") "  +0.09%  "
Set-TextValue $ws.Range("E42") "  -1.87%  "
Set-TextValue $ws.Range("D43") "0.900"
Set-TextValue $ws.Range("E43") "  -2.65%  "
Set-TextValue $ws.Range("E44") "  +2.46%  "
Set-TextValue $ws.Range("D45") "5.71"
Set-TextValue $ws.Range("E45") "  -0.12%  "
Set-TextValue $ws.Range("D46") "34.36"
Set-TextValue $ws.Range("E46") "  +31.13%  "
Set-TextValue $ws.Range("D47") "2.27"
Set-TextValue $ws.Range("E47") "  +2.18%  "
Set-TextValue $ws.Range("E48") "  -1.71%  "
Set-TextValue $ws.Range("E49") "  -5.73%  "
Set-TextValue $ws.Range("E50") "  +0.07%  "
Set-TextValue $ws.Range("D51") "8.10"
Set-TextValue $ws.Range("E51") "  -1.72%  "
